$wb = $excel.ActiveWorkbook

# --- Sheet "Instructions": update fiscal-year label in A16 ---
# (sheet is protected, so it must be unprotected before the edit and
# re-protected afterwards to keep the workbook's original state)
$wsInstructions = $wb.Worksheets.Item("Instructions")
$wsInstructions.Unprotect()
$wsInstructions.Range("A16").Value = "FY2021 Q2 (D)"
$wsInstructions.Protect()

# --- Sheet "African Development Bank": clear the D2:D11 disbursement values ---
$wsData = $wb.Worksheets.Item("African Development Bank")
$wsData.Range("D2:D11").ClearContents()
